$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 103
$ws.Range("H4").Value = 1141
$ws.Range("I4").Value = 1026
$ws.Range("J4").Value = 1108
$ws.Range("Q4").Value = 809
